# Append one new data row (date=45982) to each of the six worksheets,
# mirroring the style (date number format) of column A from the row above
# and carrying the new "remn_amt" value into column B.

$wb = $excel.ActiveWorkbook

# sheetIndex -> (newRowNumber, dateSerial, remnAmt)
$updates = @(
    @(1, 102, 45982, 700577),
    @(2, 102, 45982, 2064421),
    @(3, 102, 45982, 333245),
    @(4, 102, 45982, 149531),
    @(5, 67,  45982, 26000),
    @(6, 102, 45982, 50054)
)

foreach ($u in $updates) {
    $sheetIdx = $u[0]
    $newRow   = $u[1]
    $dateVal  = $u[2]
    $amtVal   = $u[3]
    $prevRow  = $newRow - 1

    $ws = $wb.Worksheets.Item($sheetIdx)

    # Copy the formatting (date style) of the cell above into the new cell,
    # then overwrite the values.
    $ws.Cells.Item($prevRow, 1).Copy($ws.Cells.Item($newRow, 1))

    $ws.Cells.Item($newRow, 1).Value = $dateVal
    $ws.Cells.Item($newRow, 2).Value = $amtVal
}
